# Updated cryptos list - apply per-cell edits matching the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    # Force the cell to remain plain text even when the value looks numeric
    # (prices like '1.000' or '29.203.18' must not be coerced to a Double).
    $c = $ws.Range($cellAddr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '29.197.88'
$ws.Range('E2').Value = '  -0.96%  '

Set-TextValue 'D3' '1.860.34'
$ws.Range('E3').Value = '  -1.09%  '

$ws.Range('E4').Value = '  +0.02%  '

Set-TextValue 'D5' '0.7093'
$ws.Range('E5').Value = '  -1.07%  '

Set-TextValue 'D6' '240.11'
$ws.Range('E6').Value = '  -0.85%  '

Set-TextValue 'D7' '1.001'
$ws.Range('E7').Value = '  +0.03%  '

Set-TextValue 'D8' '0.3071'
$ws.Range('E8').Value = '  -1.87%  '

Set-TextValue 'D9' '0.07631'
$ws.Range('E9').Value = '  -2.95%  '

Set-TextValue 'D10' '24.65'
$ws.Range('E10').Value = '  -1.98%  '

Set-TextValue 'D11' '0.08363'
$ws.Range('E11').Value = '  +1.21%  '

Set-TextValue 'D12' '1.862.46'
$ws.Range('E12').Value = '  -0.58%  '

Set-TextValue 'D13' '5.186'
$ws.Range('E13').Value = '  -2.01%  '

Set-TextValue 'D14' '0.7093'
$ws.Range('E14').Value = '  -3.16%  '

Set-TextValue 'D15' '90.98'
$ws.Range('E15').Value = '  -0.34%  '

Set-TextValue 'D16' '29.202.23'
$ws.Range('E16').Value = '  -1.05%  '

Set-TextValue 'D17' '5.932'
$ws.Range('E17').Value = '  -0.15%  '

Set-TextValue 'D18' '242.01'
$ws.Range('E18').Value = '  -2.59%  '

Set-TextValue 'D19' '0.000007814'
$ws.Range('E19').Value = '  -1.03%  '

Set-TextValue 'D20' '2.118.29'
$ws.Range('E20').Value = '  -0.75%  '

Set-TextValue 'D21' '13.11'
$ws.Range('E21').Value = '  -1.50%  '

$ws.Range('E22').Value = '  +0.05%  '

Set-TextValue 'D23' '7.832'
$ws.Range('E23').Value = '  -1.86%  '

$ws.Range('E24').Value = '  +0.07%  '

Set-TextValue 'D25' '0.1586'
$ws.Range('E25').Value = '  -0.64%  '

Set-TextValue 'D26' '162.58'
$ws.Range('E26').Value = '  -1.01%  '

Set-TextValue 'D27' '8.874'
$ws.Range('E27').Value = '  -1.88%  '

Set-TextValue 'D28' '18.43'
$ws.Range('E28').Value = '  +0.55%  '

Set-TextValue 'D29' '1.338'
$ws.Range('E29').Value = '  -1.86%  '

$ws.Range('E30').Value = '  +0.12%  '

Set-TextValue 'D31' '4.390'
$ws.Range('E31').Value = '  +0.29%  '

Set-TextValue 'D32' '4.199'
$ws.Range('E32').Value = '  +1.54%  '

Set-TextValue 'D33' '0.05118'
$ws.Range('E33').Value = '  -3.64%  '

Set-TextValue 'D34' '0.7992'
$ws.Range('E34').Value = '  +10.36%  '

Set-TextValue 'D35' '1.921'
$ws.Range('E35').Value = '  -0.93%  '

Set-TextValue 'D36' '1.163'
$ws.Range('E36').Value = '  -3.35%  '

Set-TextValue 'D37' '2.686'
$ws.Range('E37').Value = '  +0.29%  '

Set-TextValue 'D38' '0.01842'
$ws.Range('E38').Value = '  -1.59%  '

Set-TextValue 'D39' '2.694'
$ws.Range('E39').Value = '  -1.21%  '

Set-TextValue 'D40' '1.173.46'
$ws.Range('E40').Value = '  -7.26%  '

Set-TextValue 'D41' '6.181'
$ws.Range('E41').Value = '  +1.21%  '

Set-TextValue 'D42' '0.8893'
$ws.Range('E42').Value = '  -2.44%  '

Set-TextValue 'D43' '72.47'
$ws.Range('E43').Value = '  -2.46%  '

Set-TextValue 'D44' '1.000'
$ws.Range('E44').Value = '  -0.04%  '

Set-TextValue 'D45' '101.68'
$ws.Range('E45').Value = '  -2.03%  '

Set-TextValue 'D46' '2.015.54'
$ws.Range('E46').Value = '  -0.73%  '

$ws.Range('E47').Value = '  -2.87%  '

Set-TextValue 'D48' '1.778'
$ws.Range('E48').Value = '  +0.33%  '

$ws.Range('E49').Value = '  -0.25%  '

Set-TextValue 'D50' '9.222'
$ws.Range('E50').Value = '  -0.76%  '

$ws.Range('B51').Value = 'Frax'
$ws.Range('C51').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue 'D51' '0.9953'
$ws.Range('E51').Value = '  -0.19%  '

